# Updates cryptos list (price + 1h volume change columns) to the latest
# scraped snapshot. D-column prices that would otherwise be mis-parsed as
# numbers (losing meaningful trailing/decimal digits, e.g. "20.00" ->
# 20 or "628.40" -> 628.4) are explicitly formatted as text first, just
# like typing into a Text-formatted cell in Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "76.421.87"
$ws.Range("E2").Value = "  +0.38%  "

$ws.Range("D3").Value = "3.024.93"
$ws.Range("E3").Value = "  +3.60%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "200.13"
$ws.Range("E5").Value = "  +0.36%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "628.40"
$ws.Range("E6").Value = "  +4.66%  "

$ws.Range("E8").Value = "  +0.24%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.207"
$ws.Range("E9").Value = "  +3.27%  "

$ws.Range("D10").Value = "3.024.81"
$ws.Range("E10").Value = "  +3.67%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.435"
$ws.Range("E11").Value = "  +0.98%  "

$ws.Range("E12").Value = "  -0.42%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.12"
$ws.Range("E13").Value = "  +4.94%  "

$ws.Range("D14").Value = "3.584.88"
$ws.Range("E14").Value = "  +3.74%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.16"
$ws.Range("E15").Value = "  +5.84%  "

$ws.Range("D16").Value = "76.363.42"
$ws.Range("E16").Value = "  +0.47%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000191"
$ws.Range("E17").Value = "  -0.95%  "

$ws.Range("D18").Value = "3.017.87"
$ws.Range("E18").Value = "  +3.67%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.45"
$ws.Range("E19").Value = "  +3.21%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.05"
$ws.Range("E20").Value = "  +3.20%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "373.55"
$ws.Range("E21").Value = "  +0.55%  "

$ws.Range("E22").Value = "  -0.42%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.27"
$ws.Range("E23").Value = "  -1.49%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "72.99"
$ws.Range("E24").Value = "  +2.39%  "

$ws.Range("D25").Value = "3.187.31"
$ws.Range("E25").Value = "  +3.93%  "

$ws.Range("E26").Value = "  +0.15%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.37"
$ws.Range("E27").Value = "  +3.78%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.93"
$ws.Range("E28").Value = "  +2.61%  "

$ws.Range("E29").Value = "  -1.11%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.13%  "

$ws.Range("E31").Value = "  +7.14%  "

$ws.Range("E32").Value = "  -0.41%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "508.52"
$ws.Range("E33").Value = "  +0.94%  "

$ws.Range("E34").Value = "  +6.76%  "

$ws.Range("E35").Value = "  +0.02%  "

$ws.Range("E36").Value = "  +2.30%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "163.96"
$ws.Range("E37").Value = "  -0.81%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "193.53"
$ws.Range("E38").Value = "  +6.84%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "20.01"
$ws.Range("E39").Value = "  +1.80%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.383"
$ws.Range("E40").Value = "  +10.05%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.104"
$ws.Range("E41").Value = "  -0.20%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.112"
$ws.Range("E42").Value = "  -1.79%  "

$ws.Range("E43").Value = "  +0.37%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.07"
$ws.Range("E44").Value = "  +1.54%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.43"
$ws.Range("E45").Value = "  +5.69%  "

$ws.Range("E46").Value = "  +5.39%  "

$ws.Range("E47").Value = "  -0.03%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.712"
$ws.Range("E48").Value = "  +7.86%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.603"
$ws.Range("E49").Value = "  +5.06%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.35"
$ws.Range("E50").Value = "  +0.29%  "

$ws.Range("E51").Value = "  +3.85%  "

